$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (fill/border/font/row-height) of row 21 onto row 22,
# since row 22 moves from the "pink" table-row style to the "blue"
# supplementary-material style used by row 21.
$ws.Range("A21:E21").Copy() | Out-Null
$ws.Range("A22:E22").PasteSpecial(-4122) | Out-Null

# Update the cell contents of row 22 to the new figures/tables list entry.
$ws.Range("A22").Value = "Supplementary Table S3.3"
$ws.Range("B22").Value = "Online Supplementary Material"
$ws.Range("C22").Value = "List of all Eggnog orthogroups and ATFDB annotation."
$ws.Range("D22").Value = "Prepared. To be added in github."
$ws.Range("E22").Value = ""

# Match the row height of row 21 (default, non-wrapped) rather than the old
# taller height that was needed for the two-line text that used to be here.
$ws.Range("A22").EntireRow.AutoFit() | Out-Null

# Move the active selection, matching the author's last cursor position.
$ws.Range("B24").Select() | Out-Null
